$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.810.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "'3.164.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -7.62%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'567.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'169.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.32%  "
$ws.Range("D7").Value = "'0.618"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'3.165.77"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("E10").Value = "  -6.12%  "
$ws.Range("D11").Value = "'6.55"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.14%  "
$ws.Range("D12").Value = "'0.394"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.88%  "
$ws.Range("D13").Value = "'3.714.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.65%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'27.08"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.28%  "
$ws.Range("D16").Value = "'64.792.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("E17").Value = "  -6.36%  "
$ws.Range("D18").Value = "'3.167.70"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -7.96%  "
$ws.Range("D19").Value = "'5.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "'12.83"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.46%  "
$ws.Range("D21").Value = "'356.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").Value = "'7.27"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.51%  "
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'69.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").Value = "'0.497"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.89%  "
$ws.Range("E26").Value = "  -7.54%  "
$ws.Range("D27").Value = "'9.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.98"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.95%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'5.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.06%  "
$ws.Range("D34").Value = "'1.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.50%  "
$ws.Range("D35").Value = "'6.64"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("D36").Value = "'1.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "'154.92"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "'0.837"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("D39").Value = "'26.25"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "'2.48"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("D42").Value = "'2.659.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'4.19"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("D45").Value = "'39.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "'24.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").Value = "'323.61"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").Value = "'0.0273"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.82%  "
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.06%  "
